$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("97:97").Copy()
$ws.Rows("98:98").Insert()

$ws.Range("A98").Value = 45326.75380677084
$ws.Range("B98").Value = "Bhuvanesh Ganta "
$ws.Range("C98").Value = "B23485 "
$ws.Range("D98").Value = "https://www.beecrowd.com.br/judge/en/profile/948316"
$ws.Range("E98").Value = "VLSI"
$ws.Range("F98").Value = 0
$ws.Range("F98").ClearFormats()
$ws.Range("F98").Font.ThemeColor = 1

$ws.Hyperlinks.Add($ws.Range("D98"), "https://www.beecrowd.com.br/judge/en/profile/948316")
$ws.Range("D98").Font.Name = "Arial"
$ws.Range("D98").Font.Underline = 2
$ws.Range("D98").Font.Color = 13391121

$ws.Hyperlinks.Add($ws.Range("E98"), "https://www.beecrowd.com.br/judge/en/profile/948316")
$ws.Range("E98").Font.Name = "Arial"
$ws.Range("E98").Font.Underline = 2
$ws.Range("E98").Font.Color = 16711680

Write-Host "done"
